$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit swaps the full contents of row 4 <-> row 5, and row 15 <-> row 16.
# Only the columns that actually differ between each pair are rewritten
# (A, B, D, E, F, G, H, Q, R); columns that hold identical values in both
# rows of a pair (I, P, S, T, U, V, W, Y, AA, AD, AE, AG, AT, AW, AX, AY)
# are intentionally left untouched.

function Swap-RowValues($ws, $rowA, $rowB, $cols) {
    foreach ($col in $cols) {
        $addrA = "$col$rowA"
        $addrB = "$col$rowB"
        $valA = $ws.Range($addrA).Value()
        $valB = $ws.Range($addrB).Value()
        $ws.Range($addrA).Value = $valB
        $ws.Range($addrB).Value = $valA
    }
}

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

Swap-RowValues $ws 4 5 $cols
Swap-RowValues $ws 15 16 $cols

# Row 4/5 also differ in the presence of the sparse, empty "AF" cell
# (Bestamningsmetod): originally only AF5 carries an empty placeholder
# cell; after the edit it belongs to AF4 instead.
$ws.Range("AF5").Copy($ws.Range("AF4"))
$ws.Range("AF5").ClearContents()
